# Update capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 entirely (data for Unternehmens Invest AG and startup300 AG removed)
$ws.Rows("4:5").Delete()

# Row 2 updates
$ws.Range("B2").Formula = "'1"
$ws.Range("B2").Style = $ws.Range("C2").Style
$ws.Range("D2").Value = 0.0351
$ws.Range("E2").Value = -0.353
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.01020886293241032
$ws.Range("J2").Value = 0.006532500527769016
$ws.Range("K2").Value = 0.338
$ws.Range("L2").Value = 0.01769633507853403
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 0.083
$ws.Range("V2").Value = 0.002561728395061729
$ws.Range("W2").Value = 0.006787148594377511
$ws.Range("X2").Value = 0.07979111419601945
$ws.Range("Y2").Value = -0.07300396560164193
$ws.Range("Z2").Value = 0.1709094970334014
$ws.Range("AA2").Value = 0.001116466379571432
$ws.Range("AB2").Value = 0.03783753091140812
$ws.Range("AC2").Value = -0.03672106453183668
$ws.Range("AD2").Value = 62.3
$ws.Range("AE2").Value = 3.030053589954814
$ws.Range("AF2").Value = 65.33005358995482
$ws.Range("AG2").Value = 65.24705358995482
$ws.Range("AH2").Value = 0.6684745499481621
$ws.Range("AI2").Value = 0.5649919857481608
$ws.Range("AJ2").Value = 0.6681927533005146
$ws.Range("AK2").Value = 0.5646795098860671
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 77.77777777777777
$ws.Range("AO2").ClearContents()
$ws.Range("AP2").Value = 81.45699574276506
$ws.Range("AQ2").ClearContents()

# Row 3 updates
$ws.Range("D3").Value = 0.0351
$ws.Range("E3").Value = -0.353
$ws.Range("I3").Value = 0.01020886293241032
$ws.Range("J3").Value = 0.006532500527769016
$ws.Range("K3").Value = 0.338
$ws.Range("L3").Value = 0.01769633507853403
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 0.083
$ws.Range("V3").Value = 0.002561728395061729
$ws.Range("W3").Value = 0.006787148594377511
$ws.Range("X3").Value = 0.07979111419601945
$ws.Range("Y3").Value = -0.07300396560164193
$ws.Range("Z3").Value = 0.1709094970334014
$ws.Range("AA3").Value = 0.001116466379571432
$ws.Range("AB3").Value = 0.03783753091140812
$ws.Range("AC3").Value = -0.03672106453183668
$ws.Range("AD3").Value = 62.3
$ws.Range("AE3").Value = 3.030053589954814
$ws.Range("AF3").Value = 65.33005358995482
$ws.Range("AG3").Value = 65.24705358995482
$ws.Range("AH3").Value = 0.6684745499481621
$ws.Range("AI3").Value = 0.5649919857481608
$ws.Range("AJ3").Value = 0.6681927533005146
$ws.Range("AK3").Value = 0.5646795098860671
$ws.Range("AN3").Value = 77.77777777777777
$ws.Range("AP3").Value = 81.45699574276506
